# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 156 (pushing all existing rows 156-188
# down by one, to 157-189), then populate the new row 156 with this week's
# values. The previous last row (old 188) naturally lands at 189 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 156..188 down to 157..189.
$ws.Rows.Item(156).Insert()

# Fill in the new row 156 with the new weekly record.
$ws.Cells.Item(156, 1).Value = 5
$ws.Cells.Item(156, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(156, 3).Value = "Maule"
$ws.Cells.Item(156, 4).Value = 44504
$ws.Cells.Item(156, 5).Value = 7
$ws.Cells.Item(156, 6).Value = 100112006
$ws.Cells.Item(156, 7).Value = "Repollo"
$ws.Cells.Item(156, 8).Value = "Crespo record"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 6000
$ws.Cells.Item(156, 11).Value = 600
$ws.Cells.Item(156, 12).Value = 600
$ws.Cells.Item(156, 13).Value = 600
$ws.Cells.Item(156, 14).Value = "$/unidad"
$ws.Cells.Item(156, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(156, 16).Value = 600
$ws.Cells.Item(156, 17).Value = 1
$ws.Cells.Item(156, 18).Value = "Hortaliza"
